$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix casing typo: Claimheader -> ClaimHeader for existing rows (2-5)
$ws.Range("B2").Value = "InsuranceClaim.Mediclaim.ClaimHeader.claimType"
$ws.Range("B3").Value = "InsuranceClaim.Mediclaim.ClaimHeader.corporateEntityCode"
$ws.Range("B4").Value = "InsuranceClaim.Mediclaim.ClaimHeader.iparTransmissionRule"
$ws.Range("B5").Value = "InsuranceClaim.Mediclaim.ClaimHeader.discountType"

# Add new rows 6-11
$ws.Range("A6").Value = "(Claim) Total Allowed Amount"
$ws.Range("B6").Value = "InsuranceClaim.Mediclaim.ClaimHeader.totalAllowedAmount"

$ws.Range("A7").Value = "(Claim) Ineligible Amount"
$ws.Range("B7").Value = "InsuranceClaim.Mediclaim.ClaimHeader.ineligibleAmount"

$ws.Range("A8").Value = "(Claim) Total Billed Amount"
$ws.Range("B8").Value = "InsuranceClaim.Mediclaim.ClaimHeader.totalBilledAmount"

$ws.Range("A9").Value = "(Claim) Eligible Amount"
$ws.Range("B9").Value = "InsuranceClaim.Mediclaim.ClaimHeader.eligibleAmount"

$ws.Range("A10").Value = "(Claim) Ineligible Reason Code"
$ws.Range("B10").Value = "InsuranceClaim.Mediclaim.ClaimHeader.ineligibleReasonCode"

$ws.Range("A11").Value = "(Claim) State"
$ws.Range("B11").Value = "InsuranceClaim.Mediclaim.ClaimHeader.state"

$ws.Range("B8").Select()
